$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (order matters for shared-string table insertion order) ---
$ws.Range("O1").Value = "TotRes"
$ws.Range("T1").Value = "DistanceMeasured"
$ws.Range("S1").Value = "ScFactor"
$ws.Range("R1").Value = "OldScaledFactor"

# --- Column width for S (column 19) ---
$ws.Columns.Item(19).ColumnWidth = 11.81640625

# --- Data rows: R (raw pixel distance), S (=2.9091*R, metres), T (=1.92*R, metres) ---
$ws.Range("R2").Value = 33.854199999999999
$ws.Range("S2").Formula = "=2.9091*R2"
$ws.Range("T2").Formula = "=1.92*R2"
$ws.Range("R3").Value = 33.854199999999999
$ws.Range("S3").Formula = "=2.9091*R3"
$ws.Range("T3").Formula = "=1.92*R3"
$ws.Range("R4").Value = 33.854199999999999
$ws.Range("S4").Formula = "=2.9091*R4"
$ws.Range("T4").Formula = "=1.92*R4"
$ws.Range("R5").Value = 33.854199999999999
$ws.Range("S5").Formula = "=2.9091*R5"
$ws.Range("T5").Formula = "=1.92*R5"
$ws.Range("R6").Value = 11.452999999999999
$ws.Range("S6").Formula = "=2.9091*R6"
$ws.Range("T6").Formula = "=1.92*R6"
$ws.Range("R7").Value = 33.854199999999999
$ws.Range("S7").Formula = "=2.9091*R7"
$ws.Range("T7").Formula = "=1.92*R7"
$ws.Range("R8").Value = 33.854199999999999
$ws.Range("S8").Formula = "=2.9091*R8"
$ws.Range("T8").Formula = "=1.92*R8"
$ws.Range("R9").Value = 33.854199999999999
$ws.Range("S9").Formula = "=2.9091*R9"
$ws.Range("T9").Formula = "=1.92*R9"
$ws.Range("R10").Value = 11.452999999999999
$ws.Range("S10").Formula = "=2.9091*R10"
$ws.Range("T10").Formula = "=1.92*R10"
$ws.Range("R11").Value = 33.854199999999999
$ws.Range("S11").Formula = "=2.9091*R11"
$ws.Range("T11").Formula = "=1.92*R11"
$ws.Range("R12").Value = 23.43
$ws.Range("S12").Formula = "=2.9091*R12"
$ws.Range("T12").Formula = "=1.92*R12"
$ws.Range("R13").Value = 33.33
$ws.Range("S13").Formula = "=2.9091*R13"
$ws.Range("T13").Formula = "=1.92*R13"
$ws.Range("R14").Value = 33.33
$ws.Range("S14").Formula = "=2.9091*R14"
$ws.Range("T14").Formula = "=1.92*R14"
$ws.Range("R15").Value = 33.33
$ws.Range("S15").Formula = "=2.9091*R15"
$ws.Range("T15").Formula = "=1.92*R15"
$ws.Range("R16").Value = 33.33
$ws.Range("S16").Formula = "=2.9091*R16"
$ws.Range("T16").Formula = "=1.92*R16"
$ws.Range("R17").Value = 26.05
$ws.Range("S17").Formula = "=2.9091*R17"
$ws.Range("T17").Formula = "=1.92*R17"
$ws.Range("R18").Value = 44.27
$ws.Range("S18").Formula = "=2.9091*R18"
$ws.Range("T18").Formula = "=1.92*R18"
$ws.Range("R19").Value = 16.48
$ws.Range("S19").Formula = "=2.9091*R19"
$ws.Range("T19").Formula = "=1.92*R19"
$ws.Range("R20").Value = 27.08
$ws.Range("S20").Formula = "=2.9091*R20"
$ws.Range("T20").Formula = "=1.92*R20"
$ws.Range("R21").Value = 27.08
$ws.Range("S21").Formula = "=2.9091*R21"
$ws.Range("T21").Formula = "=1.92*R21"
$ws.Range("R22").Value = 27.08
$ws.Range("S22").Formula = "=2.9091*R22"
$ws.Range("T22").Formula = "=1.92*R22"
$ws.Range("R23").Value = 27.08
$ws.Range("S23").Formula = "=2.9091*R23"
$ws.Range("T23").Formula = "=1.92*R23"
$ws.Range("R24").Value = 27.08
$ws.Range("S24").Formula = "=2.9091*R24"
$ws.Range("T24").Formula = "=1.92*R24"
$ws.Range("R25").Value = 27.08
$ws.Range("S25").Formula = "=2.9091*R25"
$ws.Range("T25").Formula = "=1.92*R25"
$ws.Range("R26").Value = 27.08
$ws.Range("S26").Formula = "=2.9091*R26"
$ws.Range("T26").Formula = "=1.92*R26"
$ws.Range("R27").Value = 27.08
$ws.Range("S27").Formula = "=2.9091*R27"
$ws.Range("T27").Formula = "=1.92*R27"
$ws.Range("R28").Value = 27.08
$ws.Range("S28").Formula = "=2.9091*R28"
$ws.Range("T28").Formula = "=1.92*R28"
$ws.Range("R29").Value = 10.41
$ws.Range("S29").Formula = "=2.9091*R29"
$ws.Range("T29").Formula = "=1.92*R29"
$ws.Range("R30").Value = 32.82
$ws.Range("S30").Formula = "=2.9091*R30"
$ws.Range("T30").Formula = "=1.92*R30"
$ws.Range("R31").Value = 32.82
$ws.Range("S31").Formula = "=2.9091*R31"
$ws.Range("T31").Formula = "=1.92*R31"
$ws.Range("R32").Value = 11.452999999999999
$ws.Range("S32").Formula = "=2.9091*R32"
$ws.Range("T32").Formula = "=1.92*R32"
$ws.Range("R33").Value = 32.82
$ws.Range("S33").Formula = "=2.9091*R33"
$ws.Range("T33").Formula = "=1.92*R33"
$ws.Range("R34").Value = 32.82
$ws.Range("S34").Formula = "=2.9091*R34"
$ws.Range("T34").Formula = "=1.92*R34"
$ws.Range("R35").Value = 32.82
$ws.Range("S35").Formula = "=2.9091*R35"
$ws.Range("T35").Formula = "=1.92*R35"
$ws.Range("R36").Value = 6.68
$ws.Range("S36").Formula = "=2.9091*R36"
$ws.Range("T36").Formula = "=1.92*R36"
$ws.Range("R37").Value = 32.82
$ws.Range("S37").Formula = "=2.9091*R37"
$ws.Range("T37").Formula = "=1.92*R37"
$ws.Range("R38").Value = 6.68
$ws.Range("S38").Formula = "=2.9091*R38"
$ws.Range("T38").Formula = "=1.92*R38"
$ws.Range("R39").Value = 32.82
$ws.Range("S39").Formula = "=2.9091*R39"
$ws.Range("T39").Formula = "=1.92*R39"
$ws.Range("R40").Value = 32.82
$ws.Range("S40").Formula = "=2.9091*R40"
$ws.Range("T40").Formula = "=1.92*R40"
$ws.Range("R41").Value = 33.854199999999999
$ws.Range("S41").Formula = "=2.9091*R41"
$ws.Range("T41").Formula = "=1.92*R41"
$ws.Range("R42").Value = 11.452999999999999
$ws.Range("S42").Formula = "=2.9091*R42"
$ws.Range("T42").Formula = "=1.92*R42"
$ws.Range("R43").Value = 29.68
$ws.Range("S43").Formula = "=2.9091*R43"
$ws.Range("T43").Formula = "=1.92*R43"
$ws.Range("R44").Value = 29.68
$ws.Range("S44").Formula = "=2.9091*R44"
$ws.Range("T44").Formula = "=1.92*R44"
$ws.Range("R45").Value = 99.479200000000006
$ws.Range("S45").Formula = "=2.9091*R45"
$ws.Range("T45").Formula = "=1.92*R45"
$ws.Range("R46").Value = 50
$ws.Range("S46").Formula = "=2.9091*R46"
$ws.Range("T46").Formula = "=1.92*R46"
$ws.Range("R47").Value = 99.479200000000006
$ws.Range("S47").Formula = "=2.9091*R47"
$ws.Range("T47").Formula = "=1.92*R47"
$ws.Range("R48").Value = 33.854199999999999
$ws.Range("S48").Formula = "=2.9091*R48"
$ws.Range("T48").Formula = "=1.92*R48"
$ws.Range("R49").Value = 33.854199999999999
$ws.Range("S49").Formula = "=2.9091*R49"
$ws.Range("T49").Formula = "=1.92*R49"
$ws.Range("R50").Value = 33.854199999999999
$ws.Range("S50").Formula = "=2.9091*R50"
$ws.Range("T50").Formula = "=1.92*R50"
$ws.Range("R51").Value = 33.854199999999999
$ws.Range("S51").Formula = "=2.9091*R51"
$ws.Range("T51").Formula = "=1.92*R51"
$ws.Range("R52").Value = 33.854199999999999
$ws.Range("S52").Formula = "=2.9091*R52"
$ws.Range("T52").Formula = "=1.92*R52"
$ws.Range("R53").Value = 33.33
$ws.Range("S53").Formula = "=2.9091*R53"
$ws.Range("T53").Formula = "=1.92*R53"
$ws.Range("R54").Value = 27.08
$ws.Range("S54").Formula = "=2.9091*R54"
$ws.Range("T54").Formula = "=1.92*R54"
$ws.Range("R55").Value = 32.82
$ws.Range("S55").Formula = "=2.9091*R55"
$ws.Range("T55").Formula = "=1.92*R55"
$ws.Range("R56").Value = 11.452999999999999
$ws.Range("S56").Formula = "=2.9091*R56"
$ws.Range("T56").Formula = "=1.92*R56"
$ws.Range("R57").Value = 33.854199999999999
$ws.Range("S57").Formula = "=2.9091*R57"
$ws.Range("T57").Formula = "=1.92*R57"
$ws.Range("R58").Value = 33.854199999999999
$ws.Range("S58").Formula = "=2.9091*R58"
$ws.Range("T58").Formula = "=1.92*R58"
$ws.Range("R59").Value = 33.854199999999999
$ws.Range("S59").Formula = "=2.9091*R59"
$ws.Range("T59").Formula = "=1.92*R59"
$ws.Range("R60").Value = 10.41
$ws.Range("S60").Formula = "=2.9091*R60"
$ws.Range("T60").Formula = "=1.92*R60"
$ws.Range("R61").Value = 33.854199999999999
$ws.Range("S61").Formula = "=2.9091*R61"
$ws.Range("T61").Formula = "=1.92*R61"
$ws.Range("R62").Value = 33.854199999999999
$ws.Range("S62").Formula = "=2.9091*R62"
$ws.Range("T62").Formula = "=1.92*R62"
$ws.Range("R63").Value = 33.854199999999999
$ws.Range("S63").Formula = "=2.9091*R63"
$ws.Range("T63").Formula = "=1.92*R63"
$ws.Range("R64").Value = 10.414999999999999
$ws.Range("S64").Formula = "=2.9091*R64"
$ws.Range("T64").Formula = "=1.92*R64"
$ws.Range("R65").Value = 23.43
$ws.Range("S65").Formula = "=2.9091*R65"
$ws.Range("T65").Formula = "=1.92*R65"
$ws.Range("R66").Value = 23.43
$ws.Range("S66").Formula = "=2.9091*R66"
$ws.Range("T66").Formula = "=1.92*R66"
$ws.Range("R67").Value = 33.33
$ws.Range("S67").Formula = "=2.9091*R67"
$ws.Range("T67").Formula = "=1.92*R67"
$ws.Range("R68").Value = 33.33
$ws.Range("S68").Formula = "=2.9091*R68"
$ws.Range("T68").Formula = "=1.92*R68"
$ws.Range("R69").Value = 10.41
$ws.Range("S69").Formula = "=2.9091*R69"
$ws.Range("T69").Formula = "=1.92*R69"
$ws.Range("R70").Value = 20.83
$ws.Range("S70").Formula = "=2.9091*R70"
$ws.Range("T70").Formula = "=1.92*R70"
$ws.Range("R71").Value = 33.33
$ws.Range("S71").Formula = "=2.9091*R71"
$ws.Range("T71").Formula = "=1.92*R71"
$ws.Range("R72").Value = 32.82
$ws.Range("S72").Formula = "=2.9091*R72"
$ws.Range("T72").Formula = "=1.92*R72"
$ws.Range("R73").Value = 32.82
$ws.Range("S73").Formula = "=2.9091*R73"
$ws.Range("T73").Formula = "=1.92*R73"
$ws.Range("R74").Value = 32.82
$ws.Range("S74").Formula = "=2.9091*R74"
$ws.Range("T74").Formula = "=1.92*R74"
$ws.Range("R75").Value = 32.82
$ws.Range("S75").Formula = "=2.9091*R75"
$ws.Range("T75").Formula = "=1.92*R75"
$ws.Range("R76").Value = 41.6
$ws.Range("S76").Formula = "=2.9091*R76"
$ws.Range("T76").Formula = "=1.92*R76"
$ws.Range("R77").Value = 41.6
$ws.Range("S77").Formula = "=2.9091*R77"
$ws.Range("T77").Formula = "=1.92*R77"
$ws.Range("R78").Value = 41.6
$ws.Range("S78").Formula = "=2.9091*R78"
$ws.Range("T78").Formula = "=1.92*R78"

# --- Reset S:T number formatting that may have been inherited from R during formula entry ---
$ws.Range("S2:T78").Style = "Normal"

# --- Apply number format to R cells that need the custom format ---
$ws.Range("R2:R5,R7:R9,R11,R41,R48:R53,R57:R64").NumberFormat = "#,##0.0000"
$ws.Range("R45,R47").NumberFormat = "#,##0"

# --- Selection / view state ---
$ws.Range("S1").Select()